$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title heading and matching bold run near the end (both identical, replaced identically)
Replace-Text "Play Mermaid Queen for Free - Slot Game Review" "Play Mermaid Queen - Free Games and Big Wins!"

# "What we like" bullet list
Replace-Text "Free games feature that can be reactivated countless times" "Simple gameplay suitable for all players"
Replace-Text "Multiplier bonus of 6x for any win involving a wild symbol during free games" "Free games feature with a multiplier bonus"
Replace-Text "Suitable for all types of players with betting options ranging from €0.10 up to €500" "Reactivate the free games feature countless times"
Replace-Text "Welcome Bonuses and free spins available from various online casinos" "Numerous other Mermaid-themed slots available"

# "What we don't like" bullet list
Replace-Text "Outdated and less dynamic in terms of design" "Outdated and less dynamic graphics"

# Meta description (italic run)
Replace-Text "Discover Mermaid Queen, a Barcrest slot game with a free games feature and multiplier bonuses. Play for free and check out other Mermaid-themed slots." "Discover the Mermaid Queen slot game with a free games feature and multiplier bonus - play for free now!"
